$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.971.34'
$ws.Range("E2").Value = '  -4.21%  '

$ws.Range("D3").Value = '1.740.53'
$ws.Range("E3").Value = '  -4.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.54'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5796'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2737'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.19'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06628'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07543'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.74%  '

$ws.Range("D12").Value = '1.738.70'
$ws.Range("E12").Value = '  -4.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.710'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6021'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.08%  '

$ws.Range("D15").Value = '1.977.32'
$ws.Range("E15").Value = '  -4.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '74.71'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008765'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -10.84%  '

$ws.Range("D18").Value = '27.946.14'
$ws.Range("E18").Value = '  -3.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.318'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '205.55'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.30'
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.634'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.15%  '

$ws.Range("E24").Value = '  -0.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.30'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.049'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.31%  '

$ws.Range("E27").Value = '  -4.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.17'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.385'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06168'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.93%  '

$ws.Range("E31").Value = '  -3.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.746'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.738'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.06%  '

$ws.Range("E34").Value = '  -2.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.038'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6407'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.74%  '

$ws.Range("E37").Value = '  -4.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.717'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01668'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.76%  '

$ws.Range("D40").Value = '1.125.81'
$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("E41").Value = '  -6.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8740'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.92'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("D45").Value = '1.889.45'
$ws.Range("E45").Value = '  -4.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.39'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.39%  '

$ws.Range("E47").Value = '  -4.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.581'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.274'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05378'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.14%  '

$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.258'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.55%  '
